# Change the AIO courier price on the "aio" sheet and leave it as the
# active sheet/selection (the workbook previously had "polcar" active).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("aio")

$priceCell = $ws.Range("C2")
$priceCell.Value = 14.63

# The original cell used a "0.00" number format; the updated workbook
# stores it with the default "General" format (while keeping the
# existing wrap-text alignment), so reset the style accordingly.
$priceCell.Style = "Normal"
$priceCell.WrapText = $true

# Make "aio" the active sheet/tab, with the selection moved to D4 (as in
# the target workbook), and leave "polcar" on its previous selection.
[void]$ws.Range("D4").Select()
